$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("A1").Value = "Indice Viaje"

# Update data values
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 4
$ws.Range("A5").Value = 5

# Remove the old extra rows (A6:A9)
$ws.Range("A6:A9").ClearContents()

# Update selection to N13
$ws.Range("N13").Select()
